$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '23.800.98'
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -3.54%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.618.28'
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -3.52%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '307.57'
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -2.06%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.9997'
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3917'
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -0.54%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3841'
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -2.88%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.378'
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -2.05%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '49.20'
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -3.05%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08445'
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -2.51%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '23.95'
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -5.31%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '7.042'
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -4.26%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.543'
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -2.40%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.00001277'
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -3.36%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.622.59'
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -3.41%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '93.47'
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -0.61%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06910'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '20.05'
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -5.18%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.811'
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -3.93%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.9992'
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '13.41'
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  -3.93%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '23.817.50'
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -3.50%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.426'
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  +2.94%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.896'
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  +3.97%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '22.18'
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -4.04%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '156.33'
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -2.44%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '139.36'
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -5.10%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '5.282'
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -10.23%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '7.858'
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -6.28%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '2.481'
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  -0.66%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.802.44'
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -3.20%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.08078'
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -2.93%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.9834'
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -1.02%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.02888'
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -6.70%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '6.581'
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  -5.64%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.2671'
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -4.98%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.09172'
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  -3.92%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '10.36'
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '13.55'
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -0.18%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.424'
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -6.71%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.7515'
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -5.21%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '16.16'
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -2.92%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.6892'
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -3.27%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.471'
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -3.77%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '4.064'
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -2.63%  '
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.08256'
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -4.65%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '133.66'
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -3.06%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.217'
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -8.61%  '
